$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): C1/D1/E1 rotate values ---
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# --- Data rows: column C becomes the species string (was numeric duplicate of B) ---
# --- column D stays the species string (unchanged) ---
# --- column E becomes the new numeric rejection-f probability (was species string) ---

$species = "s__CAG-345 sp000433315"

$ws.Range("C2").Value = $species
$ws.Range("C3").Value = $species
$ws.Range("C4").Value = $species
$ws.Range("C5").Value = $species
$ws.Range("C6").Value = $species
$ws.Range("C7").Value = $species
$ws.Range("C8").Value = $species
$ws.Range("C9").Value = $species
$ws.Range("C10").Value = $species

$ws.Range("E2").Value = 0.9894000735264343
$ws.Range("E3").Value = 0.9889412405195864
$ws.Range("E4").Value = 0.988858088633821
$ws.Range("E5").Value = 0.9890246509343378
$ws.Range("E6").Value = 0.9893950287478965
$ws.Range("E7").Value = 0.9889412405195864
$ws.Range("E8").Value = 0.9894211488724438
$ws.Range("E9").Value = 0.989315410884634
$ws.Range("E10").Value = 0.9891086554632376

$wb.Save()
